$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 72500
$ws.Range("J3").Value = 72500
$ws.Range("L3").Value = 72500
$ws.Range("N3").Value = -72728
$ws.Range("H9").Value = 166
$ws.Range("I9").Value = 143.5
$ws.Range("J9").Value = 188.5
$ws.Range("K9").Value = 143.5
$ws.Range("L9").Value = 188.5
$ws.Range("M9").Value = 25.5
$ws.Range("N9").Value = -526.5
$ws.Range("H42").Value = 53.5
$ws.Range("I42").Value = 53.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 160.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 69.5
$ws.Range("N42").ClearContents()
$ws.Range("H55").Value = 1238.8
$ws.Range("I55").Value = 1232.6666
$ws.Range("J55").Value = 1248
$ws.Range("K55").Value = 1232.6666
$ws.Range("L55").Value = 1248
$ws.Range("M55").Value = -1018.6666
$ws.Range("N55").Value = -1676
$ws.Range("H70").Value = 3266.6667
$ws.Range("I70").Value = 3214.2856
$ws.Range("K70").Value = 9642.856800000001
$ws.Range("M70").Value = -9372.856800000001
$ws.Range("H73").Value = 3266.6667
$ws.Range("I73").Value = 3214.2856
$ws.Range("K73").Value = 9642.856800000001
$ws.Range("M73").Value = -8706.856800000001
$ws.Range("H88").Value = 3440.7
$ws.Range("J88").Value = 3489.6667
$ws.Range("L88").Value = 3489.6667
$ws.Range("N88").Value = -4301.6667
$ws.Range("H91").Value = 3440.7
$ws.Range("J91").Value = 3489.6667
$ws.Range("L91").Value = 3489.6667
$ws.Range("N91").Value = -6297.6667
$ws.Range("H102").Value = 72500
$ws.Range("J102").Value = 72500
$ws.Range("L102").Value = 72500
$ws.Range("N102").Value = -78990
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = 0

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1999.6666
$ws.Range("I21").Value = 1999.6666
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1999.6666
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1625.6666
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 36666.332
$ws.Range("J24").Value = 36666.332
$ws.Range("L24").Value = 36666.332
$ws.Range("N24").Value = -37414.332
$ws.Range("H60").Value = 35633.332
$ws.Range("I60").Value = 14900
$ws.Range("K60").Value = 14900
$ws.Range("M60").Value = -14167
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 36666.332
$ws.Range("J100").Value = 36666.332
$ws.Range("L100").Value = 36666.332
$ws.Range("N100").Value = -38830.332

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1942
$ws.Range("J16").Value = 1942
$ws.Range("L16").Value = 1942
$ws.Range("N16").Value = -2282
$ws.Range("H99").Value = 1938.6
$ws.Range("I99").Value = 1938.6
$ws.Range("K99").Value = 1938.6
$ws.Range("M99").Value = -440.5999999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 72499.75
$ws.Range("J43").Value = 72499.75
$ws.Range("L43").Value = 72499.75
$ws.Range("N43").Value = -72867.75
$ws.Range("H92").Value = 6750
$ws.Range("J92").Value = 6750
$ws.Range("L92").Value = 6750
$ws.Range("N92").Value = -11742
$ws.Range("H101").Value = 72499.75
$ws.Range("J101").Value = 72499.75
$ws.Range("L101").Value = 72499.75
$ws.Range("N101").Value = -78989.75
$ws.Range("H134").Value = 1994
$ws.Range("I134").Value = 1994
$ws.Range("K134").Value = 5982
$ws.Range("M134").Value = -3447

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2723.7778
$ws.Range("I2").Value = 3040.625
$ws.Range("J2").Value = 189
$ws.Range("K2").Value = 3040.625
$ws.Range("L2").Value = 189
$ws.Range("M2").Value = -2927.625
$ws.Range("N2").Value = -415
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("H94").Value = 58498
$ws.Range("J94").Value = 58498
$ws.Range("L94").Value = 58498
$ws.Range("N94").Value = -59850
$ws.Range("H101").Value = 22762.5
$ws.Range("J101").Value = 22762.5
$ws.Range("L101").Value = 22762.5
$ws.Range("N101").Value = -29252.5
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 341499
$ws.Range("J46").Value = 9798.799999999999
$ws.Range("L46").Value = 9798.799999999999
$ws.Range("N46").Value = -10174.8
$ws.Range("H101").Value = 66000
$ws.Range("J101").Value = 66000
$ws.Range("L101").Value = 66000
$ws.Range("N101").Value = -72490
$ws.Range("H122").Value = 3268.111
$ws.Range("I122").Value = 3364.25
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 10092.75
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -7642.75
$ws.Range("N122").Value = -12397
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -20100

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 29499.75
$ws.Range("I55").Value = 15000
$ws.Range("K55").Value = 15000
$ws.Range("M55").Value = -14723
$ws.Range("H62").Value = 4749.75
$ws.Range("I62").Value = 3999
$ws.Range("K62").Value = 3999
$ws.Range("M62").Value = -3375
$ws.Range("H64").Value = 10526
$ws.Range("J64").Value = 10526
$ws.Range("L64").Value = 10526
$ws.Range("N64").Value = -11022
$ws.Range("H65").Value = 4749.75
$ws.Range("I65").Value = 3999
$ws.Range("K65").Value = 19995
$ws.Range("M65").Value = -16875
$ws.Range("H67").Value = 10526
$ws.Range("J67").Value = 10526
$ws.Range("L67").Value = 10526
$ws.Range("N67").Value = -12242
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 182.33333
$ws.Range("I100").Value = 173.5
$ws.Range("K100").Value = 347
$ws.Range("M100").Value = 194
$ws.Range("H133").Value = 158500
$ws.Range("J133").Value = 158500
$ws.Range("L133").Value = 158500
$ws.Range("N133").Value = -168620
$ws.Range("H136").Value = 779.5833
$ws.Range("I136").Value = 595.9091
$ws.Range("K136").Value = 1787.7273
$ws.Range("M136").Value = 762.2727
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
